# Tabulated the pseudo code: indent the statements nested inside the
# if/else blocks of checkHealthMax(), and move the second "begin if"
# so it introduces the overflow branch instead of the "newHealth > MAX_HEALTH"
# condition. Also updates the active cell selection on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B23").Value = "   newHealth <= MAX_HEALTH"
$ws.Range("B24").Value = "   then health = newHealth"
$ws.Range("B25").Value = '   Print "Player is not back to full health can still play in matches"'

$ws.Range("B27").Value = "   newHealth > MAX_HEALTH"
$ws.Range("B28").Value = "   then overByMaxVal = newHealth - maxVal"
$ws.Range("B29").Value = "end if"
$ws.Range("B30").Value = "begin if"
$ws.Range("B31").Value = "   health = newHealth - overByVal"
$ws.Range("B32").Value = '   Print "Player is back to full health"'

# Update the saved selection to match the new active cell.
$ws.Range("B32").Select() | Out-Null

$wb.Save()
